# Daily attendance processing - 2026-01-10 06:43:02
# Swap the order of "Recorded By" values from "System, dnasr281@gmail.com"
# to "dnasr281@gmail.com, System" throughout the "Recorded By" column (G).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
